$d = $word.ActiveDocument

# The CRC "Collaborations" column for the five classes (Board, Ship,
# Missile, Player, Connection) had placeholder text that is being
# cleared out. Each word is the sole content of its paragraph's run,
# so a plain Find/Replace with an empty replacement removes the run
# entirely, leaving an empty paragraph (matching the diff).

$words = @("Board", "Ship", "Missile", "Player", "Connection")

foreach ($wd in $words) {
    $d.Content.Find.Execute($wd, $true, $true, $false, $false, $false, $true, 1, $false, "", 2)
}
